$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Diebold-Mariano statistics (column C) and p-values (column D)
# for rows 2-11, per the corrected computation described in the commit.

$ws.Range("C2").Value = 1.032450680271187
$ws.Range("D2").Value = 0.3091495262086186

$ws.Range("C3").Value = 0.3627071928512556
$ws.Range("D3").Value = 0.7190687281874433

$ws.Range("C4").Value = 1.380679080153042
$ws.Range("D4").Value = 0.1763871540540469

$ws.Range("C5").Value = 0.03570438652338268
$ws.Range("D5").Value = 0.9717268773527779

$ws.Range("C6").Value = -0.9075335808413791
$ws.Range("D6").Value = 0.3705124084831874

$ws.Range("C7").Value = 0.5564158324220148
$ws.Range("D7").Value = 0.5815705969139007

$ws.Range("C8").Value = -1.240896547600747
$ws.Range("D8").Value = 0.2231381600797004

$ws.Range("C9").Value = 1.304533429818878
$ws.Range("D9").Value = 0.2008135092058718

$ws.Range("C10").Value = -0.6008030400082048
$ws.Range("D10").Value = 0.5519556370800403

$ws.Range("C11").Value = -1.299761937739767
$ws.Range("D11").Value = 0.202426114791507
